# Commit: "removed errant semicolon from export/node.html"
#
# Each node in the exported casebook ends with a paragraph styled
# "Node End" whose sole run is a non-breaking space immediately
# followed by a stray ";" (an artifact of the export/node.html
# template). Strip the trailing semicolon from every such paragraph,
# leaving just the non-breaking space behind.

$d = $word.ActiveDocument
$nbsp = [char]160

foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Node End") {
        $found = $p.Range.Find.Execute($nbsp + ";", $false, $false, $false, $false, $false, `
                                        $true, 1, $false, $nbsp, 2)
    }
}
